# Applies the scheduled market-price refresh to the affected Leve Profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 173.2
$ws.Cells.Item(33, 9).Value = 149.52942
$ws.Cells.Item(33, 10).Value = 307.33334
$ws.Cells.Item(33, 11).Value = 149.52942
$ws.Cells.Item(33, 12).Value = 307.33334
$ws.Cells.Item(33, 13).Value = 79.47058000000001
$ws.Cells.Item(33, 14).Value = -765.33334

$ws.Cells.Item(132, 8).Value = 2848.3333
$ws.Cells.Item(132, 9).Value = 2519.2856
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 11).Value = 7557.8568
$ws.Cells.Item(132, 12).Value = 12000
$ws.Cells.Item(132, 13).Value = -5027.8568
$ws.Cells.Item(132, 14).Value = -17060

$ws.Cells.Item(138, 8).Value = 3390.5908

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 999.2
$ws.Cells.Item(110, 9).Value = 924
$ws.Cells.Item(110, 11).Value = 924
$ws.Cells.Item(110, 13).Value = 1121

$ws.Cells.Item(122, 8).Value = 3959
$ws.Cells.Item(122, 9).Value = 1500
$ws.Cells.Item(122, 10).Value = 5188.5
$ws.Cells.Item(122, 11).Value = 4500
$ws.Cells.Item(122, 12).Value = 15565.5
$ws.Cells.Item(122, 13).Value = -2050
$ws.Cells.Item(122, 14).Value = -20465.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1649.8667
$ws.Cells.Item(134, 9).Value = 1596.6154
$ws.Cells.Item(134, 10).Value = 1996
$ws.Cells.Item(134, 11).Value = 4789.8462
$ws.Cells.Item(134, 12).Value = 5988
$ws.Cells.Item(134, 13).Value = -2254.8462
$ws.Cells.Item(134, 14).Value = -11058

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1468.75
$ws.Cells.Item(22, 9).Value = 1683.3334
$ws.Cells.Item(22, 10).Value = 1340
$ws.Cells.Item(22, 11).Value = 1683.3334
$ws.Cells.Item(22, 12).Value = 1340
$ws.Cells.Item(22, 13).Value = -1333.3334
$ws.Cells.Item(22, 14).Value = -2040

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()

$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()

$ws.Cells.Item(98, 8).Value = 347.85715
$ws.Cells.Item(98, 9).Value = 275
$ws.Cells.Item(98, 10).Value = 377
$ws.Cells.Item(98, 11).Value = 825
$ws.Cells.Item(98, 12).Value = 1131
$ws.Cells.Item(98, 13).Value = 673
$ws.Cells.Item(98, 14).Value = -4127

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10666.333
$ws.Cells.Item(70, 9).Value = 10666.333
$ws.Cells.Item(70, 11).Value = 10666.333
$ws.Cells.Item(70, 13).Value = -10396.333

$ws.Cells.Item(73, 8).Value = 10666.333
$ws.Cells.Item(73, 9).Value = 10666.333
$ws.Cells.Item(73, 11).Value = 10666.333
$ws.Cells.Item(73, 13).Value = -9730.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3249.5
$ws.Cells.Item(40, 9).Value = 3249.5
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 3249.5
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -3113.5
$ws.Cells.Item(40, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 723.65216
$ws.Cells.Item(136, 9).Value = 704.7619
$ws.Cells.Item(136, 10).Value = 922
$ws.Cells.Item(136, 11).Value = 2114.2857
$ws.Cells.Item(136, 12).Value = 2766
$ws.Cells.Item(136, 13).Value = 435.7143000000001
$ws.Cells.Item(136, 14).Value = -7866

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124:N124").ClearContents()

$ws.Range("H125:N125").ClearContents()

$ws.Range("H127:N127").ClearContents()

$ws.Range("H128:N128").ClearContents()

$ws.Range("H129:N129").ClearContents()

$ws.Range("H130:N130").ClearContents()

$ws.Range("H131:N131").ClearContents()

$ws.Range("H132:N132").ClearContents()

$ws.Range("H133:N133").ClearContents()

$ws.Range("H134:N134").ClearContents()

$ws.Range("H135:N135").ClearContents()

$ws.Range("H136:N136").ClearContents()

$ws.Range("H137:N137").ClearContents()

$ws.Range("H138:N138").ClearContents()

$ws.Range("H139:N139").ClearContents()

$ws.Range("H140:N140").ClearContents()

$ws.Range("H141:N141").ClearContents()
